# Microsite Education Script completed
# Appends new interviewer-history rows to the AMSIN, BETA and AMS sheets,
# and fixes up the formatting/value of AMS!A28:G28 that was left
# un-styled by a previous run.

$wb = $excel.ActiveWorkbook

function Add-DataRow($ws, $row, $templateRow, $runDate, $runTime, $sprintName, $total, $pass, $fail, $timeTaken) {
    # Force column A/C..G onto the same style as the row above (General
    # number format) and column B onto the existing date/time style by
    # copying the formats down before writing values - this avoids Excel
    # auto-detecting the "yyyy-mm-dd" text in column A as a real date.
    $srcRange = $ws.Range("A" + $templateRow + ":G" + $templateRow)
    $dstRange = $ws.Range("A" + $row + ":G" + $row)
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $aCell = $ws.Range("A" + $row)
    $aCell.NumberFormat = "@"
    $aCell.Value = $runDate

    $ws.Range("B" + $row).Value = $runTime
    $ws.Range("C" + $row).Value = $sprintName
    $ws.Range("D" + $row).Value = $total
    $ws.Range("E" + $row).Value = $pass
    $ws.Range("F" + $row).Value = $fail
    $ws.Range("G" + $row).Value = $timeTaken

    # Re-stamp column A with the template row's formatting again - setting
    # NumberFormat = "@" above creates a fresh text style, so paste the
    # original formats back over the top now that the text value is safely
    # stored (PasteSpecial formats-only never touches cell content).
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------
# AMSIN sheet: rows 47-51 (sprints 165 + 166, first/second/final cycles)
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

Add-DataRow $wsAmsin 47 46 "2022-08-02" 44775.63911009259 "165_fstcycle"  75 67 8 2.84
Add-DataRow $wsAmsin 48 47 "2022-08-03" 44776.66170651621 "165_scndcycle" 75 75 0 1.86
Add-DataRow $wsAmsin 49 48 "2022-08-04" 44777.38479631944 "165_finalrun"  75 73 2 1.75
Add-DataRow $wsAmsin 50 49 "2022-08-22" 44795.6597746875  "166fstcycle"   75 75 0 1.68
Add-DataRow $wsAmsin 51 50 "2022-08-23" 44796.89887210648 "166cyclescnd"  75 74 1 1.92

# ---------------------------------------------------------------------
# BETA sheet: rows 26-27 (sprints 165 + 166 beta)
# ---------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")

Add-DataRow $wsBeta 26 25 "2022-08-04" 44777.55698798611 "165beta"  75 75 0 2.16
Add-DataRow $wsBeta 27 26 "2022-08-24" 44797.51873453704 "166_beta" 75 75 0 2.11

# ---------------------------------------------------------------------
# AMS sheet: re-style row 28 (it was missing the shared "s=5" style),
# nudge its run-time value, then append rows 29-30 (sprint 165 + 166 live)
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# Row 27 already carries the correct style - stamp it onto row 28 so
# A28/C28/D28/E28/F28/G28 pick up the shared General-format style.
$wsAms.Range("A27:G27").Copy()
$wsAms.Range("A28:G28").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$wsAms.Range("B28").Value = 44756.81275667824

Add-DataRow $wsAms 29 28 "2022-08-04" 44777.81109689815 "165_live" 75 75 0 2.14
Add-DataRow $wsAms 30 29 "2022-08-24" 44797.91305876953 "166_live" 75 75 0 2.19

Write-Host "Appended AMSIN 47-51, BETA 26-27, AMS 28 (restyled) + 29-30"
